$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Negate values for rows 6 and 7 in columns C and F
# (detection of lines being on right or left side)
$ws.Range("C6").Value = -692.7
$ws.Range("F6").Value = -761.74
$ws.Range("C7").Value = -116.42
$ws.Range("F7").Value = -185.56
